$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("B2").Value = 45861701.46629615
$ws.Range("C2").Value = 3178.151334576864
$ws.Range("D2").Value = "[5.85108711 1.04025955]"
$ws.Range("E2").Value = 28.20407700538635

$ws.Range("B3").Value = 3017.796647366936
$ws.Range("C3").Value = 2880.12358130603
$ws.Range("D3").Value = "[ 2.71314885 18.76493152]"
$ws.Range("E3").Value = 28.20407700538635

$ws.Range("B4").Value = 2821.567765789802
$ws.Range("C4").Value = 2480.582326825271
$ws.Range("D4").Value = "[ 3.16108031 18.79074714]"
$ws.Range("E4").Value = 28.20407700538635

$ws.Range("B5").Value = 2724.360496048061
$ws.Range("C5").Value = 2480.582326825271
$ws.Range("D5").Value = "[ 3.16108031 18.79074714]"
$ws.Range("E5").Value = 28.20407700538635

$ws.Range("B6").Value = 2448.458099808952
$ws.Range("C6").Value = 1920.205556846319
$ws.Range("D6").Value = "[ 4.0256312  18.76574977]"
$ws.Range("E6").Value = 28.20407700538635

$ws.Range("B7").Value = 2207.000191106167
$ws.Range("C7").Value = 1865.070504249676
$ws.Range("D7").Value = "[ 4.13078268 18.78928727]"
$ws.Range("E7").Value = 28.20407700538635

$ws.Range("B8").Value = 1883.101198075768
$ws.Range("C8").Value = 1820.389134475256
$ws.Range("D8").Value = "[ 4.22105208 18.91704751]"
$ws.Range("E8").Value = 28.20407700538635

$ws.Range("B9").Value = 1870.437056700439
$ws.Range("C9").Value = 1820.389134475256
$ws.Range("D9").Value = "[ 4.22105208 18.91704751]"
$ws.Range("E9").Value = 28.20407700538635

$ws.Range("B10").Value = 1870.437056700439
$ws.Range("C10").Value = 1820.389134475256
$ws.Range("D10").Value = "[ 4.22105208 18.91704751]"
$ws.Range("E10").Value = 28.20407700538635

$ws.Range("B11").Value = 1852.012089934923
$ws.Range("C11").Value = 1819.561961136669
$ws.Range("D11").Value = "[ 4.22105208 18.93284487]"
$ws.Range("E11").Value = 28.20407700538635
